# Trabajando las fechas y extensiones en twig
#
# A new bank-statement row is inserted at the very top of the sheet (row 1),
# pushing every existing row down by one (old row 1 -> row 2, ... old row
# 113 -> row 114). The brand-new row 1 receives the newest transaction
# (MASTERCARD payment), and only row 1 keeps the helper "H" column formula
# that renders the PHP array literal - the cells that used to carry that
# formula (old H1/H2/H3, shifted by the insert to H2/H3/H4) are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push every row down by inserting a fresh row at the top.
$ws.Rows.Item(1).Insert()

# 2) The date column (A) needs the same short-date display the rest of the
#    column already uses - newly inserted cells don't inherit it because the
#    date style lives on the row, not the column. Copy the format straight
#    from the row below instead of re-typing a format string, so it reuses
#    the existing style record instead of growing a new (equivalent) one.
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns D, F and G are plain text (account docs / amounts kept as strings
# so things like trailing spaces or leading zeros survive); make sure that
# holds for the new row before writing into it.
$ws.Range("D1").NumberFormat = "@"
$ws.Range("F1").NumberFormat = "@"
$ws.Range("G1").NumberFormat = "@"

# 3) Fill in the new transaction.
$ws.Range("A1").Value = 41719
$ws.Range("B1").Value = "13359401-MASTERCARD-RA-518114000072"
$ws.Range("C1").Value = "D"
$ws.Range("D1").Value = "0001621209"
$ws.Range("E1").Value = "SERVICIOS CENTRALES"
$ws.Range("F1").Value = "234.40  "
$ws.Range("G1").Value = "3.40"

# 4) Re-enter the PHP-array-builder helper formula in H1 (it was pushed out
#    of H1 by the insert, so put it back).
$ws.Range("H1").Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A1,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B1,""', 'mo_tipo' => '"",C1,""', 'mo_documento' => '"",D1,""', 'mo_oficina' => '"",E1,""', 'mo_monto' => "",F1,"", 'mo_saldo' => "",G1,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_borrado_logico' => false),"")"

# 5) The insert carried the old helper formula down into H2/H3/H4 (it used to
#    live in the now-shifted rows 1/2/3) - those no longer have a purpose,
#    only the new row keeps it.
$ws.Range("H2").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("H4").ClearContents()

# 6) Only H1 stays selected afterwards.
$ws.Range("H1").Select()
